$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. First paragraph: "In a previous post, we showed how ..." -> "We showed how ..."
#    The "previous post" hyperlink is removed entirely (its text is dropped).
# ---------------------------------------------------------------------------
$hPrevPost = $d.Hyperlinks.Item(1)
$hPrevPost.Range.Delete()
$null = $d.Content.Find.Execute(
    "In a , we showed", $true, $false, $false, $false, $false,
    $true, 1, $false, "We showed", 2)

# ---------------------------------------------------------------------------
# Helper: remove a hyperlink but keep its visible text as plain (non-hyperlink,
# non-colored/underlined) text that merges into the surrounding plain runs -
# mirrors what Word does when you use "Remove Hyperlink" and then clear the
# leftover character formatting.
# ---------------------------------------------------------------------------
function Unlink-HyperlinkKeepText($index) {
    $h = $d.Hyperlinks.Item($index)
    $txt = $h.Range.Text
    $h.Delete()

    # Locate the (now unwrapped, still colored/underlined) text run and insert
    # a brand-new plain-text copy immediately in front of it; the new text
    # inherits the formatting of the plain run it is typed into.
    $searchRange = $d.Content
    $null = $searchRange.Find.Execute($txt)
    $insStart = $searchRange.Start
    $insPoint = $d.Range($insStart, $insStart)
    $insPoint.InsertAfter($txt)

    # Remove the old, still-formatted (colored/underlined) copy of the text.
    $searchRange2 = $d.Range($insStart + $txt.Length, $d.Content.End)
    $null = $searchRange2.Find.Execute($txt)
    $searchRange2.Delete()
}

# ---------------------------------------------------------------------------
# 2. "the caret or stats packages" -> hyperlinks on "caret" and "stats" are
#    removed, the words stay as plain text.
# ---------------------------------------------------------------------------
Unlink-HyperlinkKeepText 2  # caret (index 2 now that "previous post" is gone)
Unlink-HyperlinkKeepText 2  # stats (index 2 again, after caret is unlinked)

# ---------------------------------------------------------------------------
# 3. Remove the closing "That's it for this post! ..." paragraph completely
#    (text, its hyperlink, and the paragraph mark).
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "That*s it for this post*") {
        $p.Range.Delete()
        break
    }
}
